# Insert a new weekly data row before the existing row 96, shifting the
# remaining rows (old 96-101) down to 97-102, then populate the new row
# with the latest week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(96).Insert()

$ws.Range("A96").Value = 9
$ws.Range("B96").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C96").Value = "Metropolitana"
$ws.Range("D96").Value = 44753
$ws.Range("E96").Value = 13
$ws.Range("F96").Value = 100114007
$ws.Range("G96").Value = "Jengibre"
$ws.Range("H96").Value = "Sin especificar"
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 520
$ws.Range("K96").Value = 17000
$ws.Range("L96").Value = 18000
$ws.Range("M96").Value = 17500
$ws.Range("N96").Value = "$/caja 13 kilos"
$ws.Range("O96").Value = "Perú"
$ws.Range("P96").Value = 1346
$ws.Range("Q96").Value = 13
$ws.Range("R96").Value = "Hortaliza"
